$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 18:05"

# Estados Unidos (row 4): Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes
$ws.Range("B4").Value = 1731795
$ws.Range("C4").Value = 6520
$ws.Range("D4").Value = 481156
$ws.Range("E4").Value = 1149739
$ws.Range("G4").Value = 328
$ws.Range("H4").Value = 100900

# India (row 13)
$ws.Range("B13").Value = 157453
$ws.Range("C13").Value = 6660
$ws.Range("D13").Value = 67167
$ws.Range("E13").Value = 85764
$ws.Range("G13").Value = 178
$ws.Range("H13").Value = 4522

# Polonia (row 38)
$ws.Range("B38").Value = 22473
$ws.Range("C38").Value = 399
$ws.Range("E38").Value = 11115
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 1028

# Chequia (row 55)
$ws.Range("B55").Value = 9069
$ws.Range("C55").Value = 19
$ws.Range("D55").Value = 6361
$ws.Range("E55").Value = 2391

# Luxemburgo (row 73)
$ws.Range("B73").Value = 4001
$ws.Range("C73").Value = 6
$ws.Range("D73").Value = 3791
$ws.Range("E73").Value = 100
